# HistoryScrapeStatus.xlsx - "Clean up and finished data gathering, on to data processing"
#
# The "History Status v2.0" column (F) had been left half-finished: rows
# 42-70 were still blank ("Bad"/red), and two rows (37/38) had values that
# had been swapped by mistake while entering data. This finishes the scrape
# status column so it mirrors the "Size (base links)" column (C) for every
# brand, and tidies up the view/selection left over from editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the two rows where the v2.0 status value was entered wrong ---
$ws.Range("F37").Value = 8
$ws.Range("F38").Value = 2

# --- Finish the remaining "Not Done" (Bad/red) rows: 42-70 (66 was already done) ---
$remainingRows = @(42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,67,68,69,70)
foreach ($r in $remainingRows) {
    $cVal = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 6).Value = $cVal
    $ws.Cells.Item($r, 6).Style = "Good"
}

# --- Center the "Key" legend total cell (E72) like the rest of the totals row ---
$ws.Range("E72").HorizontalAlignment = -4108

# --- Reset the view: scroll back to the top and move the selection ---
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 1
$ws.Range("H12").Select()
